$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.784.68'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').Value = '2.113.64'
$ws.Range('E3').Value = '  +6.72%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('E5').Value = '  +3.39%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5316'
$ws.Range('E7').Value = '  +3.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4390'
$ws.Range('E8').Value = '  +7.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09025'
$ws.Range('E9').Value = '  +7.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.09'
$ws.Range('E10').Value = '  +8.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.179'
$ws.Range('E11').Value = '  +4.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.99'
$ws.Range('E12').Value = '  +2.94%  '
$ws.Range('D13').Value = '2.109.11'
$ws.Range('E13').Value = '  +7.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.746'
$ws.Range('E14').Value = '  +4.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.805'
$ws.Range('E15').Value = '  +5.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '97.48'
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001130'
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('E19').Value = '  +2.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.11'
$ws.Range('E20').Value = '  +2.03%  '
$ws.Range('E21').Value = '  +0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.358'
$ws.Range('E22').Value = '  +4.79%  '
$ws.Range('D23').Value = '30.835.40'
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.33'
$ws.Range('E24').Value = '  +7.60%  '
$ws.Range('D25').Value = '2.355.04'
$ws.Range('E25').Value = '  +7.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.251'
$ws.Range('E26').Value = '  +2.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.76'
$ws.Range('E27').Value = '  +0.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.587'
$ws.Range('E28').Value = '  +9.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.86'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.13'
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('E31').Value = '  +2.90%  '
$ws.Range('E32').Value = '  +2.58%  '
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.022'
$ws.Range('E34').Value = '  +6.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.540'
$ws.Range('E35').Value = '  +18.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02607'
$ws.Range('E36').Value = '  +5.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.536'
$ws.Range('E37').Value = '  +3.14%  '
$ws.Range('B38').Value = 'FraxShare'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.561'
$ws.Range('E38').Value = '  +7.47%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06739'
$ws.Range('E39').Value = '  +3.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.80'
$ws.Range('E40').Value = '  +9.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2282'
$ws.Range('E41').Value = '  +5.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6865'
$ws.Range('E42').Value = '  +4.50%  '
$ws.Range('E43').Value = '  +2.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.14'
$ws.Range('E44').Value = '  +4.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6444'
$ws.Range('E45').Value = '  +5.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.231'
$ws.Range('E47').Value = '  +1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.670'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.278'
$ws.Range('E49').Value = '  +4.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.59'
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.71'
$ws.Range('E51').Value = '  -2.26%  '
